$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled/recalculated data
$ws.Range("F2").Value  = -4
$ws.Range("F3").Value  = -13
$ws.Range("F7").Value  = 4
$ws.Range("F12").Value = -4
$ws.Range("F14").Value = -4
$ws.Range("F15").Value = -8
$ws.Range("F16").Value = 3
$ws.Range("F20").Value = -4
